$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 482.14285
$ws.Range("I18").Value = 479.16666
$ws.Range("K18").Value = 479.16666
$ws.Range("M18").Value = -195.16666
$ws.Range("H38").Value = 100
$ws.Range("I38").Value = 100
$ws.Range("K38").Value = 300
$ws.Range("M38").Value = 72
$ws.Range("H39").Value = 256.36365
$ws.Range("I39").Value = 87.75
$ws.Range("J39").Value = 706
$ws.Range("K39").Value = 263.25
$ws.Range("L39").Value = 2118
$ws.Range("M39").Value = 32.75
$ws.Range("N39").Value = -2710
$ws.Range("H43").Value = 4326.2666
$ws.Range("J43").Value = 5544.909
$ws.Range("L43").Value = 5544.909
$ws.Range("N43").Value = -5682.909
$ws.Range("H112").Value = 9082208
$ws.Range("J112").Value = 4630713
$ws.Range("L112").Value = 13892139
$ws.Range("N112").Value = -13894355
$ws.Range("H116").Value = 31255614
$ws.Range("I116").Value = 50002980
$ws.Range("J116").Value = 10002
$ws.Range("K116").Value = 50002980
$ws.Range("L116").Value = 10002
$ws.Range("M116").Value = -49999538
$ws.Range("N116").Value = -16886
$ws.Range("H129").Value = 271234.3
$ws.Range("J129").Value = 304071.8
$ws.Range("L129").Value = 912215.3999999999
$ws.Range("N129").Value = -922215.3999999999
$ws.Range("H137").Value = 95691.19
$ws.Range("I137").Value = 113031.14
$ws.Range("J137").Value = 6514.2856
$ws.Range("K137").Value = 339093.42
$ws.Range("L137").Value = 19542.8568
$ws.Range("M137").Value = -336543.42
$ws.Range("N137").Value = -24642.8568
$ws.Range("H138").Value = 3804
$ws.Range("J138").Value = 4079.5781
$ws.Range("L138").Value = 12238.7343
$ws.Range("N138").Value = -22518.7343

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 60.2
$ws.Range("I5").Value = 60.875
$ws.Range("J5").Value = 57.5
$ws.Range("K5").Value = 60.875
$ws.Range("L5").Value = 57.5
$ws.Range("M5").Value = 51.125
$ws.Range("N5").Value = -281.5
$ws.Range("H9").Value = 1500
$ws.Range("I9").Value = 1500
$ws.Range("K9").Value = 1500
$ws.Range("M9").Value = -1330
$ws.Range("H20").Value = 1500
$ws.Range("I20").Value = 1500
$ws.Range("K20").Value = 1500
$ws.Range("M20").Value = -1230
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H32").Value = 16293.982
$ws.Range("I32").Value = 11146.932
$ws.Range("K32").Value = 11146.932
$ws.Range("M32").Value = -10859.932
$ws.Range("H102").Value = 1248.3334
$ws.Range("I102").Value = 1197.5
$ws.Range("K102").Value = 1197.5
$ws.Range("M102").Value = 424.5
$ws.Range("H110").Value = 1698.7693
$ws.Range("I110").Value = 1055.6
$ws.Range("K110").Value = 1055.6
$ws.Range("M110").Value = 989.4000000000001
$ws.Range("H122").Value = 1549.0435
$ws.Range("J122").Value = 1155
$ws.Range("L122").Value = 3465
$ws.Range("N122").Value = -8365
$ws.Range("H132").Value = 7938.5977
$ws.Range("I132").Value = 1530.5571
$ws.Range("J132").Value = 45318.832
$ws.Range("K132").Value = 4591.6713
$ws.Range("L132").Value = 135956.496
$ws.Range("M132").Value = -2061.6713
$ws.Range("N132").Value = -141016.496

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 60.2
$ws.Range("I4").Value = 60.875
$ws.Range("J4").Value = 57.5
$ws.Range("K4").Value = 60.875
$ws.Range("L4").Value = 57.5
$ws.Range("M4").Value = 54.125
$ws.Range("N4").Value = -287.5
$ws.Range("H22").Value = 685.73334
$ws.Range("I22").Value = 589.44446
$ws.Range("K22").Value = 589.44446
$ws.Range("M22").Value = -416.44446
$ws.Range("H105").Value = 4611925.5
$ws.Range("I105").Value = 7577250
$ws.Range("K105").Value = 7577250
$ws.Range("M105").Value = -7575503
$ws.Range("H134").Value = 3191.6206
$ws.Range("I134").Value = 2886.52
$ws.Range("K134").Value = 8659.559999999999
$ws.Range("M134").Value = -6124.559999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1104.3334
$ws.Range("J16").Value = 930
$ws.Range("L16").Value = 930
$ws.Range("N16").Value = -1504
$ws.Range("H31").Value = 7140.5454
$ws.Range("I31").Value = 4200
$ws.Range("J31").Value = 8120.727
$ws.Range("K31").Value = 4200
$ws.Range("L31").Value = 8120.727
$ws.Range("M31").Value = -3905
$ws.Range("N31").Value = -8710.726999999999
$ws.Range("H34").Value = 7140.5454
$ws.Range("I34").Value = 4200
$ws.Range("J34").Value = 8120.727
$ws.Range("K34").Value = 4200
$ws.Range("L34").Value = 8120.727
$ws.Range("M34").Value = -3998
$ws.Range("N34").Value = -8524.726999999999
$ws.Range("H57").Value = 18000
$ws.Range("J57").Value = 18000
$ws.Range("L57").Value = 18000
$ws.Range("N57").Value = -19120
$ws.Range("H58").Value = 14792.237
$ws.Range("I58").Value = 1548.7368
$ws.Range("J58").Value = 28035.736
$ws.Range("K58").Value = 1548.7368
$ws.Range("L58").Value = 28035.736
$ws.Range("M58").Value = -1345.7368
$ws.Range("N58").Value = -28441.736
$ws.Range("H86").Value = 7503.524
$ws.Range("I86").Value = 1883.4166
$ws.Range("K86").Value = 1883.4166
$ws.Range("M86").Value = -760.4166
$ws.Range("H89").Value = 7503.524
$ws.Range("I89").Value = 1883.4166
$ws.Range("K89").Value = 9417.083000000001
$ws.Range("M89").Value = -3801.083000000001
$ws.Range("H105").Value = 2297.4546
$ws.Range("I105").Value = 1221.3334
$ws.Range("K105").Value = 1221.3334
$ws.Range("M105").Value = 525.6666
$ws.Range("H110").Value = 38000
$ws.Range("J110").Value = 38000
$ws.Range("L110").Value = 38000
$ws.Range("N110").Value = -46180
$ws.Range("H113").Value = 1104.3334
$ws.Range("J113").Value = 930
$ws.Range("L113").Value = 930
$ws.Range("N113").Value = -5270
$ws.Range("H122").Value = 3118.4546
$ws.Range("I122").Value = 3602.1667
$ws.Range("J122").Value = 2538
$ws.Range("K122").Value = 10806.5001
$ws.Range("L122").Value = 7614
$ws.Range("M122").Value = -8356.500100000001
$ws.Range("N122").Value = -12514
$ws.Range("H136").Value = 14792.237
$ws.Range("I136").Value = 1548.7368
$ws.Range("J136").Value = 28035.736
$ws.Range("K136").Value = 4646.2104
$ws.Range("L136").Value = 84107.208
$ws.Range("M136").Value = -2096.2104
$ws.Range("N136").Value = -89207.208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 1000
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H131").Value = 763.04
$ws.Range("J131").Value = 763.04
$ws.Range("L131").Value = 2289.12
$ws.Range("N131").Value = -12369.12

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1899
$ws.Range("I97").Value = 1409.875
$ws.Range("J97").Value = 3855.5
$ws.Range("K97").Value = 1409.875
$ws.Range("L97").Value = 3855.5
$ws.Range("M97").Value = -913.875
$ws.Range("N97").Value = -4847.5
$ws.Range("H102").Value = 1492.2162
$ws.Range("I102").Value = 1287.5938
$ws.Range("J102").Value = 2801.8
$ws.Range("K102").Value = 1287.5938
$ws.Range("L102").Value = 2801.8
$ws.Range("M102").Value = 334.4061999999999
$ws.Range("N102").Value = -6045.8
$ws.Range("H132").Value = 69496.08
$ws.Range("I132").Value = 72604.8
$ws.Range("K132").Value = 217814.4
$ws.Range("M132").Value = -215284.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2014.4286
$ws.Range("I22").Value = 2220.2
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 2220.2
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -1925.2
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 2014.4286
$ws.Range("I27").Value = 2220.2
$ws.Range("J27").Value = 1500
$ws.Range("K27").Value = 2220.2
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = -2113.2
$ws.Range("N27").Value = -1714
$ws.Range("I132").Value = 1984.3
$ws.Range("J132").Value = 6300
$ws.Range("K132").Value = 5952.9
$ws.Range("L132").Value = 18900
$ws.Range("M132").Value = -3422.9
$ws.Range("N132").Value = -23960
$ws.Range("H136").Value = 3546.5386
$ws.Range("I136").Value = 2957.8572
$ws.Range("J136").Value = 4233.3335
$ws.Range("K136").Value = 8873.571599999999
$ws.Range("L136").Value = 12700.0005
$ws.Range("M136").Value = -6323.571599999999
$ws.Range("N136").Value = -17800.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4471.1113
$ws.Range("I62").Value = 3996.6667
$ws.Range("K62").Value = 3996.6667
$ws.Range("M62").Value = -3372.6667
$ws.Range("H65").Value = 4471.1113
$ws.Range("I65").Value = 3996.6667
$ws.Range("K65").Value = 19983.3335
$ws.Range("M65").Value = -16863.3335
$ws.Range("H122").Value = 1326.9584
$ws.Range("I122").Value = 1207
$ws.Range("J122").Value = 2166.6667
$ws.Range("K122").Value = 3621
$ws.Range("L122").Value = 6500.000100000001
$ws.Range("M122").Value = -1171
$ws.Range("N122").Value = -11400.0001
$ws.Range("H132").Value = 1373
$ws.Range("I132").Value = 982.2222
$ws.Range("J132").Value = 2779.8
$ws.Range("K132").Value = 2946.6666
$ws.Range("L132").Value = 8339.400000000001
$ws.Range("M132").Value = -416.6666
$ws.Range("N132").Value = -13399.4
